$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the S.No. column (B3:B6): 4,1,2,3 -> 1,2,3,4
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4

# Update the selected cell to C13
$ws.Range("C13").Select()
